$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Set the "Results" column (E) for rows 2-21 to "SKIP" (b suite change)
$ws.Range("E2:E21").Value = "SKIP"
